# Generate Report for Handoff
#
# The localization-status report was regenerated. Several rows that were
# still "Ready for handoff" at the time of the previous run now carry a
# refreshed Latest Handoff/Handback timestamp (rows 7, 10, 11, 12, 13, 14,
# 15, 16 on every sheet). On the "Overview" sheet that timestamp lives in
# column D ("Latest Handback DateTime"); on the "zh-cn" / "de-de" sheets it
# lives in column E ("Latest Handoff Datetime").

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# Overview sheet - column D, refreshed to 2016-03-24 09:31:14
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("D$r").Value = "2016-03-24 09:31:14"
}

# zh-cn sheet - column E, refreshed to 2016-03-24 09:31:09
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "2016-03-24 09:31:09"
}

# de-de sheet - column E, refreshed to 2016-03-24 09:31:14
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "2016-03-24 09:31:14"
}
